# Update the "Obrigatorio" (column E) flag from "N" to "S" for the
# multi-record layout fields that are now required, per the commit:
# "Implement multi-record support in layout normalization and validation".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToUpdate = @(2,3,4,5,6,7,8,9,10,11,13,20)

foreach ($r in $rowsToUpdate) {
    $ws.Range("E$r").Value = "S"
}
